# Insert a new weekly price record at row 223 (Macroferia Regional de Talca - Acelga),
# pushing the existing rows 223:274 down to 224:275.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 223:274 down by one row, leaving a blank row 223 to fill in.
$ws.Rows(223).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(223, 1).Value = 5
$ws.Cells.Item(223, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(223, 3).Value = "Maule"
$ws.Cells.Item(223, 4).Value = 44754
$ws.Cells.Item(223, 5).Value = 7
$ws.Cells.Item(223, 6).Value = 100112009
$ws.Cells.Item(223, 7).Value = "Acelga"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 400
$ws.Cells.Item(223, 11).Value = 3500
$ws.Cells.Item(223, 12).Value = 3500
$ws.Cells.Item(223, 13).Value = 3500
$ws.Cells.Item(223, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(223, 15).Value = "Región del Maule"
$ws.Cells.Item(223, 16).Value = 875
$ws.Cells.Item(223, 17).Value = 4
$ws.Cells.Item(223, 18).Value = "Hortaliza"
